$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (columns A-G)
$data = @(
    @(79432, "Caroline Barros",   "Atendimento ao Cliente", "Doença",             7, 45103, 9302.82),
    @(24355, "Sofia Rezende",     "Atendimento ao Cliente", "Consulta médica",    4, 45088, 2640.21),
    @(84387, "Augusto Ribeiro",   "Financeiro",              "Doença",             2, 45096, 3527.38),
    @(61682, "Mariana Novaes",    "Financeiro",              "Viagem de negócios", 1, 45084, 10615.58),
    @(27469, "Thiago da Mata",    "Atendimento ao Cliente", "Outros",             3, 45084, 7096.79),
    @(13364, "Ana da Mata",       "P&D",                     "Consulta médica",    2, 45105, 2508.85),
    @(82186, "Dr. Raul Gonçalves","Jurídico",                "Viagem de negócios", 4, 45091, 11774.52),
    @(86438, "Thales Ferreira",   "Recursos Humanos",        "Viagem de negócios", 2, 45104, 3661.06),
    @(85695, "Antônio Farias",    "Engenharia",              "Problemas pessoais", 2, 45098, 11570.73),
    @(96160, "Thales Cunha",      "Marketing",                "Doença",             2, 45081, 9928.639999999999)
)

$rowIndex = 2
foreach ($rowData in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowData[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowData[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowData[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rowData[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rowData[4]
    $ws.Cells.Item($rowIndex, 6).Value = $rowData[5]
    $ws.Cells.Item($rowIndex, 7).Value = $rowData[6]
    $rowIndex++
}
